$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value (Price and Volume(1h) columns updated by the crypto data refresh)
$updates = @{
    "D2" = "36.648.90"
    "E2" = "  -0.31%  "
    "D3" = "2.050.56"
    "E3" = "  -0.66%  "
    "E4" = "  +0.03%  "
    "D5" = "246.42"
    "E5" = "  +0.57%  "
    "D6" = "0.666"
    "E6" = "  +2.45%  "
    "D7" = "57.80"
    "E7" = "  +3.25%  "
    "E8" = "  +0.11%  "
    "D9" = "63.25"
    "E9" = "  +7.69%  "
    "D10" = "0.372"
    "E10" = "  +1.50%  "
    "E11" = "  -1.43%  "
    "E12" = "  -2.65%  "
    "D13" = "0.928"
    "E13" = "  +6.55%  "
    "D14" = "14.52"
    "E14" = "  -2.38%  "
    "D15" = "2.350.03"
    "E15" = "  -0.50%  "
    "D16" = "5.46"
    "E16" = "  -1.27%  "
    "D17" = "2.056.84"
    "E17" = "  -0.55%  "
    "D18" = "18.06"
    "E18" = "  +4.48%  "
    "D19" = "36.551.54"
    "E19" = "  -0.53%  "
    "D20" = "71.89"
    "E20" = "  -1.42%  "
    "D21" = "0.0₃0863"
    "E21" = "  -1.16%  "
    "D22" = "237.29"
    "E22" = "  +0.68%  "
    "D23" = "5.23"
    "E23" = "  -3.21%  "
    "E24" = "  -0.12%  "
    "D25" = "2.37"
    "E25" = "  -2.09%  "
    "D26" = "2.27"
    "E26" = "  +3.84%  "
    "D27" = "9.38"
    "E27" = "  -5.11%  "
    "D28" = "164.80"
    "E28" = "  -1.26%  "
    "D29" = "20.03"
    "E29" = "  -2.45%  "
    "E30" = "  -1.19%  "
    "D31" = "1.20"
    "E31" = "  +3.43%  "
    "D32" = "5.02"
    "E32" = "  -5.68%  "
    "E33" = "  -0.86%  "
    "D34" = "4.45"
    "E34" = "  -6.00%  "
    "D35" = "0.0879"
    "E35" = "  +3.77%  "
    "E36" = "  +0.07%  "
    "D37" = "1.82"
    "E37" = "  -0.80%  "
    "E38" = "  -5.40%  "
    "D39" = "5.14"
    "E39" = "  +4.81%  "
    "E40" = "  -4.32%  "
    "E41" = "  -1.19%  "
    "D42" = "2.89"
    "E42" = "  -1.30%  "
    "E43" = "  -3.12%  "
    "D44" = "94.42"
    "E44" = "  -1.63%  "
    "D45" = "0.0913"
    "E45" = "  -3.86%  "
    "D46" = "16.07"
    "E46" = "  -1.17%  "
    "D47" = "1.381.22"
    "E47" = "  +4.53%  "
    "D48" = "7.44"
    "E48" = "  +7.37%  "
    "E49" = "  +3.11%  "
    "E50" = "  -2.92%  "
    "D51" = "46.07"
    "E51" = "  +1.14%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "246.42") are not
    # coerced into Double values by Excel's automatic type inference, matching
    # the original inline-string cell content exactly.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Reset the style back to Normal/default so no stray number-format style
    # is left attached to the cell (keeps formatting identical to before).
    $cell.Style = "Normal"
}
